# Update the "Förändrad" (Changed) date column (C) from 45318 (2024-01-27)
# to 45319 (2024-01-28) for all data rows (2 through 27) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45318) {
        $cell.Value = 45319
    }
}
